# lr update sch sth form
# - bump form to V1.2: update settings!form_title and settings!form_id
# - w_school_id (survey row 7): change type text -> integer, add a numeric
#   range constraint + constraint message, and size the row for the
#   now-wrapped constraint message text
# - settings sheet becomes the active/selected tab instead of survey

$wb = $excel.ActiveWorkbook

# ---- survey sheet: w_school_id (row 7) ----
$survey = $wb.Worksheets.Item("survey")

$survey.Range("A7").Value = "integer"
$survey.Range("F7").Value = ". > 99 and . < 1000"
$survey.Range("G7").Value = "Must be two digit between 99 and 1000"
$survey.Rows.Item(7).RowHeight = 31.5

# selection on the survey sheet moves to the new constraint/message cells
$survey.Range("F7:G7").Select() | Out-Null

# ---- settings sheet: bump form_title / form_id to V1.2 ----
$settings = $wb.Worksheets.Item("settings")

$settings.Range("A2").Value = "(2024 Jan) - 1. SCH/STH – Site Level (School or Community) Form V1.2"
$settings.Range("B2").Value = "lr_sch_sth_impact_202401_1_school_v1_2"
$settings.Rows.Item(2).RowHeight = 31.5
$settings.Columns.Item(1).ColumnWidth = 56.25

# settings becomes the active sheet/tab, with B2 selected
$settings.Activate() | Out-Null
$settings.Range("B2").Select() | Out-Null
